{"js": "// CIV-11937 Update GA Doc template\n//\n// The \"request for information\" closing paragraph changes its instruction\n// from \"You must respond...\" to \"You should respond...\" (the rest of the\n// template text / merge-fields such as <<dateBy>>, <<judgeComments>>, etc.\n// are unchanged).\nconst body = context.document.body;\n\n// Locate the paragraph that contains the instruction so we only touch the\n// intended sentence, then replace \"must\" with \"should\" inside it.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.indexOf(\"You must respond to the request for information by\") !== -1\n);\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph to update.\");\n}\n\nconst results = target.search(\"must\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'must' in the target paragraph.\");\n}\n\nresults.items[0].insertText(\"should\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# CIV-11937 Update GA Doc template\n#\n# The \"request for information\" closing paragraph changes its instruction\n# from \"You must respond...\" to \"You should respond...\" (the rest of the\n# template text / merge-fields such as <<dateBy>>, <<judgeComments>>, etc.\n# are unchanged).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*You must respond to the request for information by*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find the target paragraph to update.\"\n}\n\n$rng = $target.Range\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"must\"\n$find.MatchWholeWord = $true\n$find.MatchCase = $true\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"should\"\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n$find.Execute($null, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, \"should\", $wdReplaceOne) | Out-Null\n"}
